$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 10 ("2.6" / last sub-item of US 2), shifting
# everything below (old rows 10-45) down by one (new rows 11-46).
$ws.Rows("10:10").Insert()

# Fill in the content of the newly inserted row 10 (new sub-item "2.7":
# "Criar tela de login"). A10/B10 stay blank (they belong to the merged
# A5:A9 / B5:B11 ranges), only C10/D10 get values.
$ws.Range("C10").Value = "Criar tela de login"

# Write "2.7" as literal text (not a number) into D10, matching how the
# other "2.x" labels in column D are stored, without leaving behind any
# new/unused cell style. We compute it as a text formula result in a
# scratch cell, copy only the resulting value into D10, then clean up.
$ws.Range("Z1").Formula = "=""2.7"""
$ws.Range("Z1").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

# Update the sheet view: drop the old pinned top-left cell and move the
# active selection from F10 to C10.
$ws.Range("C10").Select()

Write-Host "Inserted GS-2.7 row (Criar tela de login)"
